$p = $ppt.ActivePresentation

# The deck's single slide master (ppt/theme/theme1.xml, design "Integral" /
# "Red Violet" colour scheme) is being switched to the stock Office colour
# scheme ("Office Theme" / "Office"), i.e. the 12 theme colours that make up
# the master's colour scheme are replaced with the default Office palette:
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72
#
# ThemeColorScheme.Colors(n) follows the standard theme colour order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
# and writing .RGB rewrites the corresponding <a:srgbClr val="…"/> entry in
# the master's theme part (ppt/theme/theme1.xml).

$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
